$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Farmacias")

# Add a new row to the table (expands the table range + autofilter automatically)
$lo = $ws.ListObjects.Item("TablaFarmacias")
$newRow = $lo.ListRows.Add()

$newRowIndex = $newRow.Range.Row

# Fill in the new pharmacy's data
$ws.Cells.Item($newRowIndex, 1).Value = "FarmaPlus (ex-Danesa)"
$ws.Cells.Item($newRowIndex, 2).Value = "Av. Cabildo 2171"
$ws.Cells.Item($newRowIndex, 3).Value = "Belgrano"
$ws.Cells.Item($newRowIndex, 4).Value = "CABA"
$ws.Cells.Item($newRowIndex, 5).Value = "CABA"

$phoneCell = $ws.Cells.Item($newRowIndex, 8)
$phoneCell.Value = 1147873100
$phoneCell.HorizontalAlignment = -4131

$ws.Cells.Item($newRowIndex, 9).Value = "(11) 47872100"

# Move the active selection (as left by the edit) to F24
$null = $ws.Range("F24").Select()
